# The workbook is already open; grab the active workbook/sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetName1")

# Cell T2 held 40210; the author typed a new value, 50802.
$ws.Range("T2").Value = 50802

# The author's cursor/selection ended up on T2 (it had been on T3).
$ws.Activate()
$ws.Range("T2").Select()
